# Actualización automática 2025-08-13 11:15:10
#
# A new client "TIERRA GUAÑO JAIRO GABRIEL" is inserted into both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets (as row 30, pushing the
# following rows - and the totals row - down by one). Sale figures for
# "MANCHENO PINO HERVIN SANTIAGO" (row 19) are revised upward, which
# ripples into the monthly total (sheet "VENTA MENSUAL") and into the
# per-category compliance summary (sheet "CUMPLIMIENTO MENSUAL").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Revised sale figures for MANCHENO PINO HERVIN SANTIAGO (row 19)
$ws1.Cells.Item(19, 8).Value = 355.5       # H19 INODOROS
$ws1.Cells.Item(19, 9).Value = 122.4       # I19 LAVABOS
$ws1.Cells.Item(19, 13).Value = 1801.22    # M19 PORCELANATO
$ws1.Cells.Item(19, 16).Value = 180.96     # P19 NO RESURTIBLES

# Insert a new row for TIERRA GUAÑO JAIRO GABRIEL above row 30, pushing
# TOAQUIZA VILCA / VILLAFUERTE MASABANDA / ZAMBRANO CEDEÑO (and the
# totals row) down by one.
$ws1.Rows.Item(30).Insert()

$ws1.Cells.Item(30, 1).Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$ws1.Cells.Item(30, 2).Value = "TIERRA GUAÑO JAIRO GABRIEL"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(30, $col).Value = 0
}

# Totals row (was row 33, now row 34): denominator moves from 31 to 32
# clients; column P (NO RESURTIBLES) also gains a non-zero count because
# P19 went from 0 to 180.96.
$ws1.Cells.Item(34, 3).Value = "0 de 32"   # C
$ws1.Cells.Item(34, 4).Value = "3 de 32"   # D
$ws1.Cells.Item(34, 5).Value = "0 de 32"   # E
$ws1.Cells.Item(34, 6).Value = "0 de 32"   # F
$ws1.Cells.Item(34, 7).Value = "0 de 32"   # G
$ws1.Cells.Item(34, 8).Value = "1 de 32"   # H
$ws1.Cells.Item(34, 9).Value = "1 de 32"   # I
$ws1.Cells.Item(34, 10).Value = "0 de 32"  # J
$ws1.Cells.Item(34, 11).Value = "0 de 32"  # K
$ws1.Cells.Item(34, 12).Value = "2 de 32"  # L
$ws1.Cells.Item(34, 13).Value = "2 de 32"  # M
$ws1.Cells.Item(34, 14).Value = "0 de 32"  # N
$ws1.Cells.Item(34, 15).Value = "1 de 32"  # O
$ws1.Cells.Item(34, 16).Value = "1 de 32"  # P
$ws1.Cells.Item(34, 17).Value = "0 de 32"  # Q
$ws1.Cells.Item(34, 18).Value = "1 de 32"  # R

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Revised "agosto" figure for MANCHENO PINO HERVIN SANTIAGO (row 19)
$ws2.Cells.Item(19, 6).Value = 2681.8      # F19 agosto

# Same new-client row insertion as sheet 1.
$ws2.Rows.Item(30).Insert()

$ws2.Cells.Item(30, 1).Value = "ALMEIDA CUATIN JHONATHANN CARLOS"
$ws2.Cells.Item(30, 2).Value = "TIERRA GUAÑO JAIRO GABRIEL"
$ws2.Cells.Item(30, 3).Value = 0
$ws2.Cells.Item(30, 4).Value = 0
$ws2.Cells.Item(30, 5).Value = 0
$ws2.Cells.Item(30, 6).Value = 0
$ws2.Cells.Item(30, 7).Value = 0

# Totals row (was row 33, now row 34): only "agosto" changes, by the same
# amount the INODOROS + LAVABOS + PORCELANATO + NO RESURTIBLES figures
# increased for MANCHENO PINO HERVIN SANTIAGO.
$ws2.Cells.Item(34, 6).Value = 6986.02     # F34

# ---------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# INODOROS
$ws3.Cells.Item(7, 4).Value = 355.5
$ws3.Cells.Item(7, 5).Value = 1244.5
$ws3.Cells.Item(7, 6).Value = 0.2221875

# LAVABOS
$ws3.Cells.Item(8, 4).Value = 122.4
$ws3.Cells.Item(8, 5).Value = 502.6
$ws3.Cells.Item(8, 6).Value = 0.19584

# NO RESURTIBLES
$ws3.Cells.Item(10, 4).Value = 180.96
$ws3.Cells.Item(10, 5).Value = 469.29
$ws3.Cells.Item(10, 6).Value = 0.2782929642445213

# PORCELANATO
$ws3.Cells.Item(16, 4).Value = 3597.05
$ws3.Cells.Item(16, 5).Value = 18276.05
$ws3.Cells.Item(16, 6).Value = 0.1644508551599911

# TOTAL
$ws3.Cells.Item(19, 4).Value = 6986.02
$ws3.Cells.Item(19, 5).Value = 25123.26107555788
$ws3.Cells.Item(19, 6).Value = 0.2175701157419521
